# Natmi LR-pairs update (Wnt9a-Fzd10): add ECs sending/target rows per Dr Hou advice
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @("ECs", "Wnt9a", "Fzd10", "FAPs", 2, 0.6666666666666666, 0.2740863333333333, 0.822259, 0.05829365131129327, 0.05829365131129326, 1, 0.3333333333333333, 0.005445666666666667, 0.016337, 0.1561347172047327, 0.1561347172047327, 0.001492582809222222, 0.013433245283, 0.009101662762320072, 0.009101662762320068),
  @("ECs", "Wnt9a", "Fzd10", "sCs", 2, 0.6666666666666666, 0.2740863333333333, 0.822259, 0.05829365131129327, 0.05829365131129326, 3, 1, 0.02943233333333333, 0.088297, 0.8438652827952674, 0.8438652827952673, 0.008067000324777777, 0.072603002923, 0.04919198854897321, 0.04919198854897319),
  @("FAPs", "Wnt9a", "Fzd10", "FAPs", 3, 1, 3.335609, 10.006827, 0.7094291261882629, 0.7094291261882628, 1, 0.3333333333333333, 0.005445666666666667, 0.016337, 0.1561347172047327, 0.1561347172047327, 0.01816461474433333, 0.163481532699, 0.1107665159942051, 0.110766515994205),
  @("FAPs", "Wnt9a", "Fzd10", "sCs", 3, 1, 3.335609, 10.006827, 0.7094291261882629, 0.7094291261882628, 3, 1, 0.02943233333333333, 0.088297, 0.8438652827952674, 0.8438652827952673, 0.09817475595766667, 0.8835728036189999, 0.598662610194058, 0.5986626101940578),
  @("sCs", "Wnt9a", "Fzd10", "FAPs", 3, 1, 1.092126, 3.276378, 0.2322772225004438, 0.2322772225004438, 1, 0.3333333333333333, 0.005445666666666667, 0.016337, 0.1561347172047327, 0.1561347172047327, 0.005947354154000001, 0.05352618738600001, 0.03626653844820758, 0.03626653844820756),
  @("sCs", "Wnt9a", "Fzd10", "sCs", 3, 1, 1.092126, 3.276378, 0.2322772225004438, 0.2322772225004438, 3, 1, 0.02943233333333333, 0.088297, 0.8438652827952674, 0.8438652827952673, 0.032143816474, 0.289294348266, 0.1960106840522363, 0.1960106840522362)
)

for ($i = 0; $i -lt $data.Count; $i++) {
  $r = $i + 2
  $row = $data[$i]
  for ($j = 0; $j -lt $row.Count; $j++) {
    $col = [char](65 + $j)
    $ws.Range("$col$r").Value = $row[$j]
  }
}
